$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as text
# (shared string) even when the value looks like a number, matching the
# "Customer_ID" / "TC" / "PD" columns which are numeric-looking id strings
# that must remain text. We temporarily switch the cell to Text format,
# assign the value, then clear the formatting again so the cell keeps the
# default (unstyled) look used by the rest of the sheet.
function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$newRows = @(
    @("118500", "17706585", "6004"),
    @("118518", "17706586", "6020"),
    @("118498", "17706587", "1001"),
    @("118452", "17706588", "1001"),
    @("118518", "17706589", "6012"),
    @("118448", "17706590", "1047"),
    @("118518", "17706591", "1035"),
    @("118452", "17706592", "1150"),
    @("118448", "17706593", "1068"),
    @("118448", "17706594", "1005")
)

$startRow = 122
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    Set-TextValue $ws.Cells.Item($r, 1) $row[0]
    Set-TextValue $ws.Cells.Item($r, 2) $row[1]
    Set-TextValue $ws.Cells.Item($r, 3) $row[2]
}
